# Update the "Förändrad" (Changed) date column (C) from serial 45180
# (2023-09-11) to 45181 (2023-09-12) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 205 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -eq 45180) {
        $cell.Value2 = 45181
    }
}
